# Translations changes before pull from upstream.
#
# Inserts a new "table_specific_translations" worksheet right after
# "settings" (becomes the 2nd tab), populates it with the new
# string_token / text.default / text.spanish translation rows, and
# makes it the active/selected sheet (moving the tab selection away
# from "settings").

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("settings")

# Insert the new sheet immediately after "settings".
$newSheet = $wb.Worksheets.Add($null, $settings)
$newSheet.Name = "table_specific_translations"

# Header row.
$newSheet.Range("A1").Value = "string_token"
$newSheet.Range("B1").Value = "text.default"
$newSheet.Range("C1").Value = "text.spanish"

# Existing token reused from elsewhere in the workbook (is_override).
$newSheet.Range("A2").Value = "is_override"
$newSheet.Range("B2").Value = "Is Override"

$newSheet.Range("A3").Value = "click_to_deliver"
$newSheet.Range("B3").Value = "Click To Deliver"

$newSheet.Range("A4").Value = "entitlement_details"
$newSheet.Range("B4").Value = "Entitlement Details"

# Make the new sheet the active tab / selected cell, which also clears
# tabSelected from "settings".
$newSheet.Activate() | Out-Null
$newSheet.Range("D9").Select() | Out-Null
